$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 329, shifting rows 329:343 down to 330:344
$ws.Rows.Item(329).Insert()

# Fill in the new row 329 with the new weekly data entry
$ws.Cells.Item(329, 1).Value = 11
$ws.Cells.Item(329, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(329, 3).Value = "Bíobío"
$ws.Cells.Item(329, 4).Value = 44826
$ws.Cells.Item(329, 5).Value = 8
$ws.Cells.Item(329, 6).Value = 100112017
$ws.Cells.Item(329, 7).Value = "Apio"
$ws.Cells.Item(329, 8).Value = "Americana (o)"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 100
$ws.Cells.Item(329, 11).Value = 7000
$ws.Cells.Item(329, 12).Value = 7500
$ws.Cells.Item(329, 13).Value = 7250
$ws.Cells.Item(329, 14).Value = "$/docena de matas"
$ws.Cells.Item(329, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(329, 16).Value = 1208
$ws.Cells.Item(329, 17).Value = 6
$ws.Cells.Item(329, 18).Value = "Hortaliza"
